$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 882.9167
$ws.Range("I11").Value = 882.9167
$ws.Range("K11").Value = 882.9167
$ws.Range("M11").Value = -742.9167
$ws.Range("H39").Value = 2574.5
$ws.Range("I39").Value = 1371
$ws.Range("J39").Value = 10999
$ws.Range("K39").Value = 4113
$ws.Range("L39").Value = 32997
$ws.Range("M39").Value = -3817
$ws.Range("N39").Value = -33589
$ws.Range("H40").Value = 11117461
$ws.Range("I40").Value = 3916.3333
$ws.Range("K40").Value = 3916.3333
$ws.Range("M40").Value = -3741.3333
$ws.Range("H62").Value = 5581
$ws.Range("I62").Value = 2497.2
$ws.Range("J62").Value = 21000
$ws.Range("K62").Value = 2497.2
$ws.Range("L62").Value = 21000
$ws.Range("M62").Value = -1873.2
$ws.Range("N62").Value = -22248
$ws.Range("H65").Value = 5581
$ws.Range("I65").Value = 2497.2
$ws.Range("J65").Value = 21000
$ws.Range("K65").Value = 12486
$ws.Range("L65").Value = 105000
$ws.Range("M65").Value = -9366
$ws.Range("N65").Value = -111240
$ws.Range("H74").Value = 47626588
$ws.Range("I74").Value = 47626588
$ws.Range("K74").Value = 47626588
$ws.Range("M74").Value = -47625652
$ws.Range("H76").Value = 3886.75
$ws.Range("I76").Value = 3518.8
$ws.Range("J76").Value = 4500
$ws.Range("K76").Value = 3518.8
$ws.Range("L76").Value = 4500
$ws.Range("M76").Value = -3203.8
$ws.Range("N76").Value = -5130
$ws.Range("H77").Value = 47626588
$ws.Range("I77").Value = 47626588
$ws.Range("K77").Value = 238132940
$ws.Range("M77").Value = -238128260
$ws.Range("H79").Value = 3886.75
$ws.Range("I79").Value = 3518.8
$ws.Range("J79").Value = 4500
$ws.Range("K79").Value = 3518.8
$ws.Range("L79").Value = 4500
$ws.Range("M79").Value = -2426.8
$ws.Range("N79").Value = -6684
$ws.Range("H101").Value = 3865.3333
$ws.Range("J101").Value = 5649
$ws.Range("L101").Value = 16947
$ws.Range("N101").Value = -20191
$ws.Range("H132").Value = 2406.158
$ws.Range("I132").Value = 2542.1765
$ws.Range("K132").Value = 7626.529500000001
$ws.Range("M132").Value = -5096.529500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1811.2241
$ws.Range("I32").Value = 1807.9122
$ws.Range("K32").Value = 1807.9122
$ws.Range("M32").Value = -1520.9122
$ws.Range("H61").Value = 47621676
$ws.Range("I61").Value = 50002560
$ws.Range("K61").Value = 50002560
$ws.Range("M61").Value = -50002348
$ws.Range("H97").Value = 501.5
$ws.Range("I97").Value = 462.14285
$ws.Range("K97").Value = 462.14285
$ws.Range("M97").Value = 33.85714999999999
$ws.Range("H132").Value = 4046319
$ws.Range("I132").Value = 2860800.2
$ws.Range("K132").Value = 8582400.600000001
$ws.Range("M132").Value = -8579870.600000001
$ws.Range("H136").Value = 47621676
$ws.Range("I136").Value = 50002560
$ws.Range("K136").Value = 150007680
$ws.Range("M136").Value = -150005130

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3431.3635
$ws.Range("I86").Value = 3493.25
$ws.Range("K86").Value = 3493.25
$ws.Range("M86").Value = -2370.25
$ws.Range("H89").Value = 3431.3635
$ws.Range("I89").Value = 3493.25
$ws.Range("K89").Value = 17466.25
$ws.Range("M89").Value = -11850.25
$ws.Range("H99").Value = 2164.5
$ws.Range("J99").Value = 2496.5
$ws.Range("L99").Value = 2496.5
$ws.Range("N99").Value = -5492.5
$ws.Range("H134").Value = 13925607
$ws.Range("I134").Value = 13925607
$ws.Range("K134").Value = 41776821
$ws.Range("M134").Value = -41774286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3430.6365
$ws.Range("I31").Value = 2327.0408
$ws.Range("K31").Value = 2327.0408
$ws.Range("M31").Value = -2032.0408
$ws.Range("H34").Value = 3430.6365
$ws.Range("I34").Value = 2327.0408
$ws.Range("K34").Value = 2327.0408
$ws.Range("M34").Value = -2125.0408
$ws.Range("H58").Value = 38511150
$ws.Range("J58").Value = 54623.332
$ws.Range("L58").Value = 54623.332
$ws.Range("N58").Value = -55029.332
$ws.Range("H136").Value = 38511150
$ws.Range("J136").Value = 54623.332
$ws.Range("L136").Value = 163869.996
$ws.Range("N136").Value = -168969.996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1001815.8
$ws.Range("I4").Value = 1334393.2
$ws.Range("K4").Value = 4003179.6
$ws.Range("M4").Value = -4003067.6
$ws.Range("H5").Value = 67735.8
$ws.Range("J5").Value = 1629.8
$ws.Range("L5").Value = 4889.4
$ws.Range("N5").Value = -5113.4
$ws.Range("H39").Value = 2869.875
$ws.Range("J39").Value = 7200
$ws.Range("L39").Value = 21600
$ws.Range("N39").Value = -22188
$ws.Range("H131").Value = 2396.7144
$ws.Range("J131").Value = 2685.1428
$ws.Range("L131").Value = 8055.428400000001
$ws.Range("N131").Value = -18135.4284
$ws.Range("H134").Value = 1125.6154
$ws.Range("I134").Value = 1125.6154
$ws.Range("K134").Value = 3376.8462
$ws.Range("M134").Value = 1693.1538
$ws.Range("H135").Value = 67735.8
$ws.Range("J135").Value = 1629.8
$ws.Range("L135").Value = 14668.2
$ws.Range("N135").Value = -19738.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 54799
$ws.Range("J26").Value = 54799
$ws.Range("L26").Value = 54799
$ws.Range("N26").Value = -55359
$ws.Range("H50").Value = 54799
$ws.Range("J50").Value = 54799
$ws.Range("L50").Value = 54799
$ws.Range("N50").Value = -55795
$ws.Range("H125").Value = 49995
$ws.Range("J125").Value = 49995
$ws.Range("L125").Value = 49995
$ws.Range("N125").Value = -54915
$ws.Range("H132").Value = 2671854
$ws.Range("I132").Value = 2852826.8
$ws.Range("J132").Value = 17590
$ws.Range("K132").Value = 8558480.399999999
$ws.Range("L132").Value = 52770
$ws.Range("M132").Value = -8555950.399999999
$ws.Range("N132").Value = -57830

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2560.6072
$ws.Range("I16").Value = 1280
$ws.Range("K16").Value = 1280
$ws.Range("M16").Value = -1110
$ws.Range("H93").Value = 850.5862
$ws.Range("I93").Value = 849.7692
$ws.Range("K93").Value = 849.7692
$ws.Range("M93").Value = 398.2308
$ws.Range("H104").Value = 15226.375
$ws.Range("J104").Value = 15226.375
$ws.Range("L104").Value = 15226.375
$ws.Range("N104").Value = -22214.375
$ws.Range("H132").Value = 8069837
$ws.Range("I132").Value = 8338681.5
$ws.Range("K132").Value = 25016044.5
$ws.Range("M132").Value = -25013514.5
$ws.Range("H133").Value = 99998.336
$ws.Range("J133").Value = 99998.336
$ws.Range("L133").Value = 99998.336
$ws.Range("N133").Value = -105058.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6274.25
$ws.Range("I62").Value = 3066.3333
$ws.Range("J62").Value = 8199
$ws.Range("K62").Value = 3066.3333
$ws.Range("L62").Value = 8199
$ws.Range("M62").Value = -2442.3333
$ws.Range("N62").Value = -9447
$ws.Range("H65").Value = 6274.25
$ws.Range("I65").Value = 3066.3333
$ws.Range("J65").Value = 8199
$ws.Range("K65").Value = 15331.6665
$ws.Range("L65").Value = 40995
$ws.Range("M65").Value = -12211.6665
$ws.Range("N65").Value = -47235
$ws.Range("H81").Value = 1299.5
$ws.Range("I81").Value = 1299.5
$ws.Range("K81").Value = 2599
$ws.Range("M81").Value = -1538
$ws.Range("H84").Value = 1299.5
$ws.Range("I84").Value = 1299.5
$ws.Range("K84").Value = 12995
$ws.Range("M84").Value = -7691
$ws.Range("H110").Value = 258499.5
$ws.Range("J110").Value = 258499.5
$ws.Range("L110").Value = 258499.5
$ws.Range("N110").Value = -266679.5
$ws.Range("H132").Value = 16131564
$ws.Range("I132").Value = 20835618
$ws.Range("K132").Value = 62506854
$ws.Range("M132").Value = -62504324
$ws.Range("H136").Value = 20835298
$ws.Range("I136").Value = 21741028
$ws.Range("K136").Value = 65223084
$ws.Range("M136").Value = -65220534
